# Update master to output generated at c986bee
# Replace the two-digit / one-digit division practice problems in the
# single table on the page. Each of the 25 populated cells gets its
# text updated to a new equation. Cell-scoped Find/Replace (one match
# at a time, via wdReplaceOne) is used rather than a document-wide
# replace because a couple of the original equations are duplicated
# verbatim elsewhere in the table, and each occurrence must be mapped
# to its own distinct replacement.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$r = $t.Cell(1,1).Range
$r.Find.Execute("81÷9=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷7=11, 5", 1) | Out-Null
$r = $t.Cell(1,2).Range
$r.Find.Execute("70÷6=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "91÷8=11, 3", 1) | Out-Null
$r = $t.Cell(1,3).Range
$r.Find.Execute("33÷9=3, 6", $true, $false, $false, $false, $false, $true, 1, $false, "37÷4=9, 1", 1) | Out-Null
$r = $t.Cell(1,4).Range
$r.Find.Execute("16÷8=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "85÷5=17, 0", 1) | Out-Null
$r = $t.Cell(1,5).Range
$r.Find.Execute("47÷9=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "52÷9=5, 7", 1) | Out-Null
$r = $t.Cell(5,1).Range
$r.Find.Execute("46÷4=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "45÷8=5, 5", 1) | Out-Null
$r = $t.Cell(5,2).Range
$r.Find.Execute("37÷8=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=12, 3", 1) | Out-Null
$r = $t.Cell(5,3).Range
$r.Find.Execute("37÷5=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "38÷4=9, 2", 1) | Out-Null
$r = $t.Cell(5,4).Range
$r.Find.Execute("18÷8=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "15÷2=7, 1", 1) | Out-Null
$r = $t.Cell(5,5).Range
$r.Find.Execute("85÷7=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "88÷9=9, 7", 1) | Out-Null
$r = $t.Cell(9,1).Range
$r.Find.Execute("90÷9=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "69÷3=23, 0", 1) | Out-Null
$r = $t.Cell(9,2).Range
$r.Find.Execute("83÷8=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "77÷3=25, 2", 1) | Out-Null
$r = $t.Cell(9,3).Range
$r.Find.Execute("18÷9=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "13÷7=1, 6", 1) | Out-Null
$r = $t.Cell(9,4).Range
$r.Find.Execute("15÷8=1, 7", $true, $false, $false, $false, $false, $true, 1, $false, "66÷6=11, 0", 1) | Out-Null
$r = $t.Cell(9,5).Range
$r.Find.Execute("98÷4=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "73÷3=24, 1", 1) | Out-Null
$r = $t.Cell(13,1).Range
$r.Find.Execute("91÷3=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=4, 3", 1) | Out-Null
$r = $t.Cell(13,2).Range
$r.Find.Execute("30÷6=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "57÷9=6, 3", 1) | Out-Null
$r = $t.Cell(13,3).Range
$r.Find.Execute("81÷9=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "40÷9=4, 4", 1) | Out-Null
$r = $t.Cell(13,4).Range
$r.Find.Execute("63÷8=7, 7", $true, $false, $false, $false, $false, $true, 1, $false, "44÷6=7, 2", 1) | Out-Null
$r = $t.Cell(13,5).Range
$r.Find.Execute("21÷6=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "88÷9=9, 7", 1) | Out-Null
$r = $t.Cell(17,1).Range
$r.Find.Execute("41÷8=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "67÷7=9, 4", 1) | Out-Null
$r = $t.Cell(17,2).Range
$r.Find.Execute("43÷6=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "37÷3=12, 1", 1) | Out-Null
$r = $t.Cell(17,3).Range
$r.Find.Execute("46÷8=5, 6", $true, $false, $false, $false, $false, $true, 1, $false, "50÷5=10, 0", 1) | Out-Null
$r = $t.Cell(17,4).Range
$r.Find.Execute("86÷4=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "43÷8=5, 3", 1) | Out-Null
$r = $t.Cell(17,5).Range
$r.Find.Execute("81÷5=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "99÷8=12, 3", 1) | Out-Null

Write-Host "Done updating table cells."
